$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row10 = @("Jasmine Matchawate","-1","-2","1","-2","1","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","-2","-2","-2","-2","-1","-1","-1","-1")
$row11 = @("Noah Williams","0","0","-1","-1","0","-2","-2","0","0","-1","-1","0","0","-1","-1","0","0","-1","-1","0","0","-1","-1","0","0","-1","-1","0","0","-1","-1")

for ($c = 0; $c -lt $row10.Length; $c++) {
    if ($c -eq 0) {
        $ws.Cells.Item(10, $c + 1).Value = $row10[$c]
    } else {
        $ws.Cells.Item(10, $c + 1).Value = "'" + $row10[$c]
    }
}

for ($c = 0; $c -lt $row11.Length; $c++) {
    if ($c -eq 0) {
        $ws.Cells.Item(11, $c + 1).Value = $row11[$c]
    } else {
        $ws.Cells.Item(11, $c + 1).Value = "'" + $row11[$c]
    }
}

# Strip the auto-applied "quote prefix" style from the numeric-looking text
# cells so they keep the default (no explicit style) formatting, matching
# the rest of the data rows.
$ws.Range("B10:AF11").Style = "Normal"
